$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per the latest cryptos snapshot.
# Numeric-looking Price strings must be forced back to text (they are
# locale-formatted strings like "26.336.99", not real numbers) so Excel
# does not silently reinterpret/round them as floating point numbers.

$ws.Range("D2").Value = '26.336.99'
$ws.Range("E2").Value = '  +0.92%  '

$ws.Range("D3").Value = '1.680.04'
$ws.Range("E3").Value = '  +0.77%  '

$ws.Range("E4").Value = '  +0.21%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.19'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.69%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5279'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.31%  '

$ws.Range("E7").Value = '  +0.21%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2696'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.48%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06472'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.90%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.97'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.21%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07510'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.13%  '

$ws.Range("D12").Value = '1.695.07'
$ws.Range("E12").Value = '  +1.64%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.516'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.15%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5792'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.40%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000008520'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.54%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.81'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.70%  '

$ws.Range("D17").Value = '26.342.34'
$ws.Range("E17").Value = '  +0.69%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.929'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.24%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.88'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.90%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '189.90'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.51%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.209'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.07%  '

$ws.Range("E23").Value = '  +0.16%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '144.90'
$ws.Range("D24").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '7.792'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.09%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1256'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.43%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.79'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.14%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06527'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.55%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.364'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.40%  '

$ws.Range("E30").Value = '  +0.40%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.592'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.93%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.589'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.15%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.661'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.50%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.029'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.00%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6220'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.42%  '

$ws.Range("E36").Value = '  +1.83%  '

$ws.Range("E37").Value = '  +1.79%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.263'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.05%  '

$ws.Range("D39").Value = '1.116.13'
$ws.Range("E39").Value = '  +3.75%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01623'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.53%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8741'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.52%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.016'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.63%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.54'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.22%  '

$ws.Range("D44").Value = '1.828.97'

$ws.Range("E45").Value = '  -1.41%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '56.93'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.30%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.155'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.20%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.006'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.13%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05269'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.17%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4293'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.03%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.084'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.31%  '
